$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns keep their existing text formatting so
# numeric-looking strings (e.g. "298.29") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Update cryptocurrency price/volume data per upstream source refresh
$ws.Range("D2").Value = "45.500.17"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "2.325.95"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "298.29"
$ws.Range("E5").Value = "  -2.40%  "
$ws.Range("D6").Value = "96.55"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  +0.26%  "
$ws.Range("D9").Value = "0.505"
$ws.Range("E9").Value = "  -5.64%  "
$ws.Range("D10").Value = "33.89"
$ws.Range("E10").Value = "  -5.15%  "
$ws.Range("D11").Value = "0.0779"
$ws.Range("E11").Value = "  -4.02%  "
$ws.Range("D12").Value = "6.99"
$ws.Range("E12").Value = "  -5.75%  "
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("D14").Value = "2.697.96"
$ws.Range("E14").Value = "  +0.14%  "
$ws.Range("D15").Value = "2.374.91"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "13.51"
$ws.Range("E16").Value = "  -4.38%  "
$ws.Range("D17").Value = "0.797"
$ws.Range("E17").Value = "  -4.26%  "
$ws.Range("D18").Value = "45.527.97"
$ws.Range("E18").Value = "  -2.92%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0959"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "12.33"
$ws.Range("E20").Value = "  -9.14%  "
$ws.Range("D21").Value = "5.84"
$ws.Range("E21").Value = "  -5.09%  "
$ws.Range("D22").Value = "65.38"
$ws.Range("E22").Value = "  -3.07%  "
$ws.Range("D23").Value = "241.27"
$ws.Range("E23").Value = "  -4.09%  "
$ws.Range("D24").Value = "2.77"
$ws.Range("E24").Value = "  -6.76%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").Value = "1.85"
$ws.Range("E26").Value = "  -6.67%  "
$ws.Range("D27").Value = "39.82"
$ws.Range("E27").Value = "  -6.45%  "
$ws.Range("E28").Value = "  -1.98%  "
$ws.Range("D29").Value = "9.53"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").Value = "20.22"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").Value = "3.50"
$ws.Range("E31").Value = "  +10.64%  "
$ws.Range("D32").Value = "2.74"
$ws.Range("E32").Value = "  +4.31%  "
$ws.Range("D33").Value = "144.09"
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "5.29"
$ws.Range("E34").Value = "  -8.41%  "
$ws.Range("D35").Value = "0.0762"
$ws.Range("E35").Value = "  -6.60%  "
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -4.78%  "
$ws.Range("E37").Value = "  -3.70%  "
$ws.Range("D38").Value = "1.73"
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D39").Value = "14.99"
$ws.Range("E39").Value = "  +7.77%  "
$ws.Range("D40").Value = "3.80"
$ws.Range("E40").Value = "  -5.01%  "
$ws.Range("D41").Value = "0.0293"
$ws.Range("E41").Value = "  -6.32%  "
$ws.Range("D42").Value = "3.08"
$ws.Range("E42").Value = "  -8.96%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("D44").Value = "1.845.57"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("D45").Value = "91.32"
$ws.Range("E45").Value = "  +2.46%  "
$ws.Range("D46").Value = "1.79"
$ws.Range("E46").Value = "  -8.85%  "
$ws.Range("D47").Value = "0.181"
$ws.Range("E47").Value = "  -7.38%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "68.86"
$ws.Range("E48").Value = "  -8.04%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.571.92"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "7.89"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "94.49"
$ws.Range("E51").Value = "  -4.47%  "
